# "add technical summary to word template"
#
# The CV template's "Professional experience" table (the first table in
# the document) ends with a blank stub row (two empty cells) left over
# for a Technical Summary entry. Fill it in:
#   - left cell  -> italic label "Technical Summary"
#   - right cell -> bold-marked placeholder field "[Work_TS]"
# matching the look (fonts/lang) of the sibling rows (Project/Position/
# Responsibilities -> [Work_Project]/[Work_Position]/[Work_Resp]).
#
# Note: this runtime's Tables.Item() acts like a single shared cursor -
# grabbing Tables.Item(N) for a different N invalidates earlier saved
# references to other indices - so we resolve the table exactly once
# and then only ever navigate from that one object.

$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

# Locate the trailing empty row (both cells blank) - the stub row meant
# to hold the technical summary. Falls back to the last row if every
# row already has content.
$targetRow = 0
for ($r = 1; $r -le $table.Rows.Count; $r++) {
    $c1 = $table.Cell($r, 1).Range.Text.Trim([char]7, [char]13)
    $c2 = $table.Cell($r, 2).Range.Text.Trim([char]7, [char]13)
    if ($c1 -eq "" -and $c2 -eq "") {
        $targetRow = $r
    }
}
if ($targetRow -eq 0) {
    $targetRow = $table.Rows.Count
}

$cell1 = $table.Cell($targetRow, 1)
$cell2 = $table.Cell($targetRow, 2)

# Replace each (empty) cell paragraph with one carrying the right pPr/rPr
# and a run with the new text, via a minimal WordprocessingML package -
# InsertXML replaces the exact range's contents with the supplied markup.

$xml1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:i/><w:lang w:val="en-US" w:eastAsia="ja-JP"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:i/><w:lang w:val="en-US" w:eastAsia="ja-JP"/></w:rPr><w:t>Technical Summary</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$xml2 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:right="-1"/><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:lang w:eastAsia="ja-JP"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:val="en-US" w:eastAsia="ja-JP"/></w:rPr><w:t>[Work_TS]</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$cell1.Range.InsertXML($xml1)
$cell2.Range.InsertXML($xml2)
